$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for "volume" and "area" columns, matching the
# bold/bordered header style already used by B1:E1.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F1").Value = "volume"
$ws.Range("G1").Value = "area"

# Daily reservoir volume and area values for rows 2-114.
$volumeValues = @(1188.982978408122,1175.44415719192,1166.708590436828,1158.090844496391,1156.907481516007,1153.584709694915,1144.712194893297,1133.618369250074,1122.50410534907,1113.360897478052,1101.153147688985,1093.961773122769,1089.30846611544,1080.966017677954,1072.10571016804,1063.185669755128,1060.886892141496,1053.240732939694,1048.138026882144,1065.186606334921,1059.269174342304,1052.256069268239,1047.171241507474,1045.073634084376,1044.346956839151,1039.985579862514,1034.889929685316,1028.395681203122,1020.56664103615,1013.587519732172,1003.993418111583,994.2216873814868,988.871277634922,982.4060378324644,977.9844970645275,969.8272871927287,962.8918133718394,955.8415696015326,950.5923102491031,943.131468668358,935.7440406765031,924.6763822762965,916.3649669884915,908.4143490476893,898.0893797159897,887.4666670413496,879.7104552210012,869.4095324345225,861.1160071422103,851.3641666846715,841.8550685442772,833.2793902590663,824.2100051674004,817.2397223050752,806.148375511284,795.3353397097869,784.7372345883505,776.9195207310263,771.4170631712277,764.1052669340095,755.9965357178241,746.2850697711997,736.7713545244999,730.0710370400648,729.2452675209138,719.6187125736603,710.9662051818111,699.9895841454261,689.9205164459847,678.60542847938,666.7212012088191,656.0679969837217,644.9949142504249,632.9400100268635,623.9780759459215,617.1857662882237,608.7844373149278,599.7906709775898,591.5563422672756,582.5736948852229,575.3106288821155,566.2785786032071,556.8284315621086,549.4613599060788,540.9128761604936,532.2504298992185,522.5970772711821,511.7456351095011,504.2373969179578,497.0245767957631,491.5691799685308,484.8393635204643,476.1248597848783,467.0599907487915,456.927651618905,446.5616595820416,441.4821716660231,432.6214769955504,421.6135141576117,414.8693138935287,406.6232078979683,398.4399490889241,393.5847884399393,387.8042379229914,381.7392589333906,371.1819731852331,364.2445449108384,355.3579495885228,349.6300508315032,343.5542266040791,339.4054141388258,335.2580517113504,327.9462455707715)
$areaValues = @(1309.12227472864,1302.02440134432,1297.418220581114,1292.853549788126,1292.225130706825,1290.458495333775,1285.725974145315,1279.777112510194,1273.781642356856,1268.822264948495,1262.161765668174,1258.217056013558,1255.656107130555,1251.048024104083,1246.12996618172,1241.153485851195,1239.866844546268,1235.574849439349,1232.699870829974,1242.272031578425,1238.960368421241,1235.020737777558,1232.154190068008,1230.969171400322,1230.55830113888,1228.088617767103,1225.195017026487,1221.49439685293,1217.013850215973,1213.001667507444,1207.457987959265,1201.777558813651,1198.652481567177,1194.862065829762,1192.26078743884,1187.442241207605,1183.325233510847,1179.12082926065,1175.977659436694,1171.491205705453,1167.026586467924,1160.295510300243,1155.206732704065,1150.311053296346,1143.91199015669,1137.278614428801,1132.402674663302,1125.883628145929,1120.598322180063,1114.340940737251,1108.193936483685,1102.611090811128,1096.665515316294,1092.066626531147,1084.69491724609,1077.443142386995,1070.271612665216,1064.940045842071,1061.165937603194,1056.122862302097,1050.49220263259,1043.695064863226,1036.978491006107,1032.213063375065,1031.623721599772,1024.719963614177,1018.461424893149,1010.447399663549,1003.020986531227,994.5876194472261,985.6270170887689,977.5022695003981,968.9621157311974,959.5511918039604,952.4760148108249,947.0678044321913,940.32254997324,933.0316548062136,926.2914939669849,918.8663684690828,912.8061763590008,905.1978915894283,897.1498576126,890.8122419627967,883.3868195870932,875.7822706461966,867.2103783743514,857.4485610989701,850.6140449933367,843.9851297977805,838.9293507042304,832.6416991193355,824.414645027481,815.7524006378677,805.9405090792048,795.7565004828799,790.7109125076831,781.8201213178528,770.6125967977788,763.6548353151156,755.0506995155957,746.4040652928935,741.2219014326341,735.0003814597179,728.4112596636207,716.7876007437216,709.0400306623906,698.9849960078913,692.4241968056475,685.3948366837832,680.5526826337482,675.6773672580345,666.9956662298696)

for ($i = 0; $i -lt $volumeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $volumeValues[$i]
    $ws.Cells.Item($row, 7).Value = $areaValues[$i]
}
